# Insert a new data row at row 128 (pushing existing rows 128..231 down to
# 129..232) and populate it with the new record, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a whole new row at position 128; this shifts all rows below it
# (and their formatting) down by one, exactly like the target diff shows.
$ws.Rows("128:128").Insert()

# Fill in the new row's values.
$ws.Cells.Item(128, 1).Value2 = 1
$ws.Cells.Item(128, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(128, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(128, 4).Value2 = 44587
$ws.Cells.Item(128, 5).Value2 = 15
$ws.Cells.Item(128, 6).Value2 = 100114013
$ws.Cells.Item(128, 7).Value2 = "Zanahoria"
$ws.Cells.Item(128, 8).Value2 = "Sin especificar"
$ws.Cells.Item(128, 9).Value2 = "Primera"
$ws.Cells.Item(128, 10).Value2 = 70
$ws.Cells.Item(128, 11).Value2 = 22000
$ws.Cells.Item(128, 12).Value2 = 23000
$ws.Cells.Item(128, 13).Value2 = 22500
$ws.Cells.Item(128, 14).Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(128, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(128, 16).Value2 = 900
$ws.Cells.Item(128, 17).Value2 = 25
$ws.Cells.Item(128, 18).Value2 = "Hortaliza"
